# Update the "Means" sheet (sheet1) - Total Cancer Risk (row 9) and
# Total Respiratory hazard quotient (row 10) with the newest airtoxics NATA data.
$wb = $excel.ActiveWorkbook

$wsMeans = $wb.Worksheets.Item("Means")
$wsMeans.Range("B9").Value = 26
$wsMeans.Range("C9").Value = 27
$wsMeans.Range("D9").Value = 30
$wsMeans.Range("E9").Value = 30
$wsMeans.Range("F9").Value = 30
$wsMeans.Range("G9").Value = 31

$wsMeans.Range("B10").Value = 0.31
$wsMeans.Range("C10").Value = 0.34
$wsMeans.Range("D10").Value = 0.4
$wsMeans.Range("E10").Value = 0.4
$wsMeans.Range("F10").Value = 0.4
$wsMeans.Range("G10").Value = 0.41

# Update the "Standard Deviations" sheet (sheet2) - same rows.
$wsSD = $wb.Worksheets.Item("Standard Deviations")
$wsSD.Range("B9").Value = 8.3
$wsSD.Range("C9").Value = 7.2
$wsSD.Range("D9").Value = 0
$wsSD.Range("E9").Value = 0
$wsSD.Range("F9").Value = 2
$wsSD.Range("G9").Value = 3.5

$wsSD.Range("B10").Value = 0.11
$wsSD.Range("C10").Value = 0.1
$wsSD.Range("D10").Value = 0.000000000000000031
$wsSD.Range("E10").Value = 0.000000000000000019
$wsSD.Range("F10").Value = 0.016
$wsSD.Range("G10").Value = 0.035
